$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 9) below the existing table.
# Write column B before column A so the shared-strings table records
# "Swag Labs" ahead of "Title=", matching the source edit.
$ws.Range("B9").Value = "Swag Labs"
$ws.Range("A9").Value = "Title="

# Carry the same left-aligned cell style used by the rest of the table
# (row 8) down onto the new row, including the still-empty C9 cell.
$ws.Range("A8:C8").Copy()
$ws.Range("A9:C9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C9").ClearContents()

# Grow the table / AutoFilter range so it covers the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C9"))

# Leave the same cell selected that was active after the edit.
$ws.Range("B9").Select()
